$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = -6
    4  = 1
    5  = -1
    6  = -2
    7  = 4
    8  = -2
    10 = -1
    11 = -2
    12 = 1
    13 = 2
    14 = -5
    16 = 3
    17 = -4
    18 = 13
    19 = 2
    20 = -5
    21 = 2
    22 = 2
    24 = 3
    25 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
